$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1424505885193952
$ws.Cells.Item(2, 4).Value = 0.02013275399991699
$ws.Cells.Item(2, 5).Value = 0.1706397611198867
$ws.Cells.Item(2, 6).Value = 0.4551076624629928
$ws.Cells.Item(2, 7).Value = 0.3005287660485791
$ws.Cells.Item(2, 8).Value = 0.4641317007753187
$ws.Cells.Item(2, 11).Value = 0.4796567021739122
$ws.Cells.Item(2, 13).Value = 0.2371909127198393
$ws.Cells.Item(2, 14).Value = 1.278496167539124
$ws.Cells.Item(2, 15).Value = 1.456480937381144
$ws.Cells.Item(3, 2).Value = 0.1329280882201687
$ws.Cells.Item(3, 4).Value = 0.01783193070128419
$ws.Cells.Item(3, 5).Value = 0.1621401311927571
$ws.Cells.Item(3, 6).Value = 0.4509585701085257
$ws.Cells.Item(3, 7).Value = 0.2977446921923459
$ws.Cells.Item(3, 8).Value = 0.4661137010476608
$ws.Cells.Item(3, 11).Value = 0.4231847918051983
$ws.Cells.Item(3, 13).Value = 0.209650902747839
$ws.Cells.Item(3, 14).Value = 1.291360057554538
$ws.Cells.Item(3, 15).Value = 1.454486899445428
$ws.Cells.Item(4, 2).Value = 0.1271497787846414
$ws.Cells.Item(4, 4).Value = 0.01641102853855614
$ws.Cells.Item(4, 5).Value = 0.157059877418952
$ws.Cells.Item(4, 6).Value = 0.4487240944589672
$ws.Cells.Item(4, 7).Value = 0.2962779892827996
$ws.Cells.Item(4, 8).Value = 0.4675549582580203
$ws.Cells.Item(4, 11).Value = 0.3883267550192784
$ws.Cells.Item(4, 13).Value = 0.1927685460267412
$ws.Cells.Item(4, 14).Value = 1.299740925498085
$ws.Cells.Item(4, 15).Value = 1.454238123305601
$ws.Cells.Item(5, 2).Value = 0.1248124873505674
$ws.Cells.Item(5, 4).Value = 0.01582997461528635
$ws.Cells.Item(5, 5).Value = 0.155024263641053
$ws.Cells.Item(5, 6).Value = 0.4478922514645802
$ws.Cells.Item(5, 7).Value = 0.2957412438193856
$ws.Cells.Item(5, 8).Value = 0.4681987170948716
$ws.Cells.Item(5, 11).Value = 0.3740764014192735
$ws.Cells.Item(5, 13).Value = 0.1858957995432675
$ws.Cells.Item(5, 14).Value = 1.30327755394697
$ws.Cells.Item(5, 15).Value = 1.454381934813512
$ws.Cells.Item(6, 2).Value = 0.1244254391136081
$ws.Cells.Item(6, 4).Value = 0.01573336976330353
$ws.Cells.Item(6, 5).Value = 0.1546883370175323
$ws.Cells.Item(6, 6).Value = 0.4477588792038176
$ws.Cells.Item(6, 7).Value = 0.2956557961521753
$ws.Cells.Item(6, 8).Value = 0.4683090222443624
$ws.Cells.Item(6, 11).Value = 0.3717074231467734
$ws.Cells.Item(6, 13).Value = 0.1847550074545197
$ws.Cells.Item(6, 14).Value = 1.30387213976023
$ws.Cells.Item(6, 15).Value = 1.454420619802022
$ws.Cells.Item(7, 2).Value = 0.1271181864972135
$ws.Cells.Item(7, 4).Value = 0.01640320039524568
$ws.Cells.Item(7, 5).Value = 0.1570322844776371
$ws.Cells.Item(7, 6).Value = 0.4487125571882444
$ws.Cells.Item(7, 7).Value = 0.2962705038882731
$ws.Cells.Item(7, 8).Value = 0.4675634116700706
$ws.Cells.Item(7, 11).Value = 0.3881347525310161
$ws.Cells.Item(7, 13).Value = 0.192675829512126
$ws.Cells.Item(7, 14).Value = 1.299788130261298
$ws.Cells.Item(7, 15).Value = 1.454239070202277
$ws.Cells.Item(8, 2).Value = 0.1391531092357638
$ws.Cells.Item(8, 4).Value = 0.0193411491858555
$ws.Cells.Item(8, 5).Value = 0.1676802031045668
$ws.Cells.Item(8, 6).Value = 0.4536120701563604
$ws.Cells.Item(8, 7).Value = 0.2995183838928668
$ws.Cells.Item(8, 8).Value = 0.46476855967191
$ws.Cells.Item(8, 11).Value = 0.4602238885101428
$ws.Cells.Item(8, 13).Value = 0.2276894871486235
$ws.Cells.Item(8, 14).Value = 1.282831551654471
$ws.Cells.Item(8, 15).Value = 1.455590832738139
$ws.Cells.Item(9, 2).Value = 0.1632907710933011
$ws.Cells.Item(9, 4).Value = 0.02503625265241993
$ws.Cells.Item(9, 5).Value = 0.1896707669518207
$ws.Cells.Item(9, 6).Value = 0.4657055809154329
$ws.Cells.Item(9, 7).Value = 0.3078182432949035
$ws.Cells.Item(9, 8).Value = 0.461066466643274
$ws.Cells.Item(9, 11).Value = 0.6000985992871222
$ws.Cells.Item(9, 13).Value = 0.2965678356974095
$ws.Cells.Item(9, 14).Value = 1.253404516404839
$ws.Cells.Item(9, 15).Value = 1.465990501885727
$ws.Cells.Item(10, 2).Value = 0.1813451557723766
$ws.Cells.Item(10, 4).Value = 0.02917878712145239
$ws.Cells.Item(10, 5).Value = 0.2065204593568595
$ws.Cells.Item(10, 6).Value = 0.4761096913153864
$ws.Cells.Item(10, 7).Value = 0.3151010831835634
$ws.Cells.Item(10, 8).Value = 0.4594297606819424
$ws.Cells.Item(10, 11).Value = 0.7019221860259108
$ws.Cells.Item(10, 13).Value = 0.3473098107271895
$ws.Cells.Item(10, 14).Value = 1.234112015567177
$ws.Cells.Item(10, 15).Value = 1.478369916169726
$ws.Cells.Item(11, 2).Value = 0.1896267521109962
$ws.Cells.Item(11, 4).Value = 0.03105403249586658
$ws.Cells.Item(11, 5).Value = 0.21433999679072
$ws.Cells.Item(11, 6).Value = 0.4811735288528922
$ws.Cells.Item(11, 7).Value = 0.318673287859113
$ws.Cells.Item(11, 8).Value = 0.4589202034386233
$ws.Cells.Item(11, 11).Value = 0.7480333347168084
$ws.Cells.Item(11, 13).Value = 0.3704247621222478
$ws.Cells.Item(11, 14).Value = 1.225839920796176
$ws.Cells.Item(11, 15).Value = 1.485034019371227
$ws.Cells.Item(12, 2).Value = 0.1927724759862457
$ws.Cells.Item(12, 4).Value = 0.03176278364555429
$ws.Cells.Item(12, 5).Value = 0.2173235275171734
$ws.Cells.Item(12, 6).Value = 0.4831386971345566
$ws.Cells.Item(12, 7).Value = 0.3200633771113672
$ws.Cells.Item(12, 8).Value = 0.4587610164365685
$ws.Cells.Item(12, 11).Value = 0.7654636287890071
$ws.Cells.Item(12, 13).Value = 0.3791824093956677
$ws.Cells.Item(12, 14).Value = 1.222779948840888
$ws.Cells.Item(12, 15).Value = 1.487706238956235
$ws.Cells.Item(13, 2).Value = 0.1920945614355247
$ws.Cells.Item(13, 4).Value = 0.03161020266514214
$ws.Cells.Item(13, 5).Value = 0.2166799691875525
$ws.Cells.Item(13, 6).Value = 0.4827133456558741
$ws.Cells.Item(13, 7).Value = 0.3197623327764916
$ws.Cells.Item(13, 8).Value = 0.4587937985341597
$ws.Cells.Item(13, 11).Value = 0.7617110969115117
$ws.Cells.Item(13, 13).Value = 0.3772960942284556
$ws.Cells.Item(13, 14).Value = 1.22343574458381
$ws.Cells.Item(13, 15).Value = 1.487124113828145
$ws.Cells.Item(14, 2).Value = 0.1898853601737756
$ws.Cells.Item(14, 4).Value = 0.0311123694294011
$ws.Cells.Item(14, 5).Value = 0.2145850026335694
$ws.Cells.Item(14, 6).Value = 0.4813342503334326
$ws.Cells.Item(14, 7).Value = 0.3187869017639287
$ws.Cells.Item(14, 8).Value = 0.4589064303148973
$ws.Cells.Item(14, 11).Value = 0.7494679622027149
$ws.Cells.Item(14, 13).Value = 0.3711451691410872
$ws.Cells.Item(14, 14).Value = 1.225586722267813
$ws.Cells.Item(14, 15).Value = 1.485250883955104
$ws.Cells.Item(15, 2).Value = 0.1885334131827534
$ws.Cells.Item(15, 4).Value = 0.0308072534731636
$ws.Cells.Item(15, 5).Value = 0.2133047042295999
$ws.Cells.Item(15, 6).Value = 0.4804957142626094
$ws.Cells.Item(15, 7).Value = 0.3181942920915048
$ws.Cells.Item(15, 8).Value = 0.4589798179614917
$ws.Cells.Item(15, 11).Value = 0.7419646235756261
$ws.Cells.Item(15, 13).Value = 0.3673781348780025
$ws.Cells.Item(15, 14).Value = 1.226913698492929
$ws.Cells.Item(15, 15).Value = 1.484122842733001
$ws.Cells.Item(16, 2).Value = 0.1808052969037988
$ws.Cells.Item(16, 4).Value = 0.02905604671610007
$ws.Cells.Item(16, 5).Value = 0.2060125625673805
$ws.Cells.Item(16, 6).Value = 0.4757854183416086
$ws.Cells.Item(16, 7).Value = 0.3148728556005409
$ws.Cells.Item(16, 8).Value = 0.4594677883226126
$ws.Cells.Item(16, 11).Value = 0.6989044284384818
$ws.Cells.Item(16, 13).Value = 0.3457998323468274
$ws.Cells.Item(16, 14).Value = 1.234662760764792
$ws.Cells.Item(16, 15).Value = 1.477955195848409
$ws.Cells.Item(17, 2).Value = 0.1760817695779622
$ws.Cells.Item(17, 4).Value = 0.02797935066843138
$ws.Cells.Item(17, 5).Value = 0.2015788053180927
$ws.Cells.Item(17, 6).Value = 0.4729805788897679
$ws.Cells.Item(17, 7).Value = 0.3129017215606922
$ws.Cells.Item(17, 8).Value = 0.4598273144317062
$ws.Cells.Item(17, 11).Value = 0.6724342056627393
$ws.Cells.Item(17, 13).Value = 0.3325704052182914
$ws.Cells.Item(17, 14).Value = 1.239545681341305
$ws.Cells.Item(17, 15).Value = 1.47443615493026
$ws.Cells.Item(18, 2).Value = 0.1733713871586247
$ws.Cells.Item(18, 4).Value = 0.02735919859427582
$ws.Cells.Item(18, 5).Value = 0.199043157150264
$ws.Cells.Item(18, 6).Value = 0.4713984598915602
$ws.Cells.Item(18, 7).Value = 0.3117923690442979
$ws.Cells.Item(18, 8).Value = 0.4600562237494188
$ws.Cells.Item(18, 11).Value = 0.657189626960303
$ws.Cells.Item(18, 13).Value = 0.3249642234812598
$ws.Cells.Item(18, 14).Value = 1.242401662325435
$ws.Cells.Item(18, 15).Value = 1.472509282825058
$ws.Cells.Item(19, 2).Value = 0.1724548146932676
$ws.Cells.Item(19, 4).Value = 0.0271490783689714
$ws.Cells.Item(19, 5).Value = 0.1981871193552607
$ws.Cells.Item(19, 6).Value = 0.4708681313234919
$ws.Cells.Item(19, 7).Value = 0.3114209475799612
$ws.Cells.Item(19, 8).Value = 0.4601375281458928
$ws.Cells.Item(19, 11).Value = 0.6520247362301461
$ws.Cells.Item(19, 13).Value = 0.322389425592732
$ws.Cells.Item(19, 14).Value = 1.243376798977806
$ws.Cells.Item(19, 15).Value = 1.471873563654498
$ws.Cells.Item(20, 2).Value = 0.1765839293872489
$ws.Cells.Item(20, 4).Value = 0.02809405660210729
$ws.Cells.Item(20, 5).Value = 0.2020492807671701
$ws.Cells.Item(20, 6).Value = 0.4732759347904576
$ws.Cells.Item(20, 7).Value = 0.3131090269744448
$ws.Cells.Item(20, 8).Value = 0.4597867531309703
$ws.Cells.Item(20, 11).Value = 0.6752540422101845
$ws.Cells.Item(20, 13).Value = 0.3339783870170834
$ws.Cells.Item(20, 14).Value = 1.239020974834425
$ws.Cells.Item(20, 15).Value = 1.474800703059429
$ws.Cells.Item(21, 2).Value = 0.1905339954528955
$ws.Cells.Item(21, 4).Value = 0.03125863243303684
$ws.Cells.Item(21, 5).Value = 0.2151997342779808
$ws.Cells.Item(21, 6).Value = 0.4817380320417115
$ws.Cells.Item(21, 7).Value = 0.3190723944201181
$ws.Cells.Item(21, 8).Value = 0.4588724312457089
$ws.Cells.Item(21, 11).Value = 0.7530649165132672
$ws.Cells.Item(21, 13).Value = 0.3729517232323403
$ws.Cells.Item(21, 14).Value = 1.224952960643037
$ws.Cells.Item(21, 15).Value = 1.485797060776662
$ws.Cells.Item(22, 2).Value = 0.19970740219145
$ws.Cells.Item(22, 4).Value = 0.03331889975186897
$ws.Cells.Item(22, 5).Value = 0.2239252247863135
$ws.Cells.Item(22, 6).Value = 0.4875459725455329
$ws.Cells.Item(22, 7).Value = 0.3231876828778297
$ws.Cells.Item(22, 8).Value = 0.4584717046704014
$ws.Cells.Item(22, 11).Value = 0.803737702664904
$ws.Cells.Item(22, 13).Value = 0.398449327437902
$ws.Cells.Item(22, 14).Value = 1.216181230757478
$ws.Cells.Item(22, 15).Value = 1.493850448647351
$ws.Cells.Item(23, 2).Value = 0.1948063022515925
$ws.Cells.Item(23, 4).Value = 0.03222003795813322
$ws.Cells.Item(23, 5).Value = 0.2192562153149638
$ws.Cells.Item(23, 6).Value = 0.4844207732011583
$ws.Cells.Item(23, 7).Value = 0.3209713070418445
$ws.Cells.Item(23, 8).Value = 0.458667576042501
$ws.Cells.Item(23, 11).Value = 0.7767095709727982
$ws.Cells.Item(23, 13).Value = 0.3848384084135859
$ws.Cells.Item(23, 14).Value = 1.220824208067356
$ws.Cells.Item(23, 15).Value = 1.489472849791554
$ws.Cells.Item(24, 2).Value = 0.1763568866120124
$ws.Cells.Item(24, 4).Value = 0.02804220161416993
$ws.Cells.Item(24, 5).Value = 0.2018365371641409
$ws.Cells.Item(24, 6).Value = 0.4731423096309442
$ws.Cells.Item(24, 7).Value = 0.3130152298240461
$ws.Cells.Item(24, 8).Value = 0.4598050217009018
$ws.Cells.Item(24, 11).Value = 0.6739792765647508
$ws.Cells.Item(24, 13).Value = 0.3333418397109824
$ws.Cells.Item(24, 14).Value = 1.239258042869984
$ws.Cells.Item(24, 15).Value = 1.474635590983354
$ws.Cells.Item(25, 2).Value = 0.15670410908713
$ws.Cells.Item(25, 4).Value = 0.02350279780696241
$ws.Cells.Item(25, 5).Value = 0.1836011912987914
$ws.Cells.Item(25, 6).Value = 0.4621674775223639
$ws.Cells.Item(25, 7).Value = 0.3053653904038498
$ws.Cells.Item(25, 8).Value = 0.4618776794692963
$ws.Cells.Item(25, 11).Value = 0.5624219016708309
$ws.Cells.Item(25, 13).Value = 0.2779107045710916
$ws.Cells.Item(25, 14).Value = 1.26095628068056
$ws.Cells.Item(25, 15).Value = 1.462346003505957

Write-Output "Done updating cells"
